# Auto-generated edit script applying the Behemoth_Profits value updates
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4060
$ws.Range("I40").Value = 3516.6667
$ws.Range("K40").Value = 3516.6667
$ws.Range("M40").Value = -3341.6667
$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H76").Value = 6763.1816
$ws.Range("J76").Value = 7985
$ws.Range("L76").Value = 7985
$ws.Range("N76").Value = -8615
$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H79").Value = 6763.1816
$ws.Range("J79").Value = 7985
$ws.Range("L79").Value = 7985
$ws.Range("N79").Value = -10169
$ws.Range("H98").Value = 142929650
$ws.Range("I98").Value = 142929650
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 142929650
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -142928152
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 142929650
$ws.Range("I122").Value = 142929650
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 428788950
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -428786500
$ws.Range("N122").ClearContents()
$ws.Range("H137").Value = 7598.4165
$ws.Range("I137").Value = 2962.1428
$ws.Range("K137").Value = 8886.428400000001
$ws.Range("M137").Value = -6336.428400000001
$ws.Range("H138").Value = 3093.9275
$ws.Range("I138").Value = 1200
$ws.Range("J138").Value = 3150.4626
$ws.Range("K138").Value = 3600
$ws.Range("L138").Value = 9451.3878
$ws.Range("M138").Value = 1540
$ws.Range("N138").Value = -19731.3878

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 788
$ws.Range("I4").Value = 519
$ws.Range("J4").Value = 1124.25
$ws.Range("K4").Value = 519
$ws.Range("L4").Value = 1124.25
$ws.Range("M4").Value = -403
$ws.Range("N4").Value = -1356.25
$ws.Range("H5").Value = 521.6667
$ws.Range("I5").Value = 531
$ws.Range("J5").Value = 517
$ws.Range("K5").Value = 531
$ws.Range("L5").Value = 517
$ws.Range("M5").Value = -419
$ws.Range("N5").Value = -741
$ws.Range("H32").Value = 10002203
$ws.Range("I32").Value = 10639705
$ws.Range("J32").Value = 14670.667
$ws.Range("K32").Value = 10639705
$ws.Range("L32").Value = 14670.667
$ws.Range("M32").Value = -10639418
$ws.Range("N32").Value = -15244.667
$ws.Range("H45").Value = 2547.7334
$ws.Range("I45").Value = 1870.3
$ws.Range("K45").Value = 1870.3
$ws.Range("M45").Value = -1493.3
$ws.Range("H112").Value = 33598
$ws.Range("J112").Value = 33598
$ws.Range("L112").Value = 33598
$ws.Range("N112").Value = -36552
$ws.Range("H132").Value = 4294.925
$ws.Range("I132").Value = 1450.7778
$ws.Range("J132").Value = 10202
$ws.Range("K132").Value = 4352.3334
$ws.Range("L132").Value = 30606
$ws.Range("M132").Value = -1822.3334
$ws.Range("N132").Value = -35666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 521.6667
$ws.Range("I4").Value = 531
$ws.Range("J4").Value = 517
$ws.Range("K4").Value = 531
$ws.Range("L4").Value = 517
$ws.Range("M4").Value = -416
$ws.Range("N4").Value = -747
$ws.Range("H20").Value = 3132.3157
$ws.Range("J20").Value = 2498
$ws.Range("L20").Value = 2498
$ws.Range("N20").Value = -2992
$ws.Range("H22").Value = 333.75
$ws.Range("I22").Value = 333.75
$ws.Range("K22").Value = 333.75
$ws.Range("M22").Value = -160.75
$ws.Range("H94").Value = 1008.1923
$ws.Range("J94").Value = 707.1111
$ws.Range("L94").Value = 707.1111
$ws.Range("N94").Value = -1609.1111
$ws.Range("H134").Value = 89189.664
$ws.Range("I134").Value = 1138.6666
$ws.Range("K134").Value = 3415.9998
$ws.Range("M134").Value = -880.9998000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 381617.9
$ws.Range("I31").Value = 4730.488
$ws.Range("J31").Value = 1117445.8
$ws.Range("K31").Value = 4730.488
$ws.Range("L31").Value = 1117445.8
$ws.Range("M31").Value = -4435.488
$ws.Range("N31").Value = -1118035.8
$ws.Range("H34").Value = 381617.9
$ws.Range("I34").Value = 4730.488
$ws.Range("J34").Value = 1117445.8
$ws.Range("K34").Value = 4730.488
$ws.Range("L34").Value = 1117445.8
$ws.Range("M34").Value = -4528.488
$ws.Range("N34").Value = -1117849.8
$ws.Range("H58").Value = 1198.2858
$ws.Range("I58").Value = 997.8
$ws.Range("J58").Value = 1699.5
$ws.Range("K58").Value = 997.8
$ws.Range("L58").Value = 1699.5
$ws.Range("M58").Value = -794.8
$ws.Range("N58").Value = -2105.5
$ws.Range("H105").Value = 2298.3
$ws.Range("I105").Value = 2109.2222
$ws.Range("K105").Value = 2109.2222
$ws.Range("M105").Value = -362.2222000000002
$ws.Range("H136").Value = 1198.2858
$ws.Range("I136").Value = 997.8
$ws.Range("J136").Value = 1699.5
$ws.Range("K136").Value = 2993.4
$ws.Range("L136").Value = 5098.5
$ws.Range("M136").Value = -443.3999999999996
$ws.Range("N136").Value = -10198.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2702.0454
$ws.Range("I11").Value = 2749.762
$ws.Range("K11").Value = 8249.286
$ws.Range("M11").Value = -8109.286
$ws.Range("H76").Value = 4988.8
$ws.Range("J76").Value = 4988.8
$ws.Range("L76").Value = 14966.4
$ws.Range("N76").Value = -15732.4
$ws.Range("H79").Value = 4988.8
$ws.Range("J79").Value = 4988.8
$ws.Range("L79").Value = 14966.4
$ws.Range("N79").Value = -17618.4
$ws.Range("H80").Value = 5014.9
$ws.Range("I80").Value = 3499.5
$ws.Range("K80").Value = 10498.5
$ws.Range("M80").Value = -9562.5
$ws.Range("H83").Value = 5014.9
$ws.Range("I83").Value = 3499.5
$ws.Range("K83").Value = 31495.5
$ws.Range("M83").Value = -26815.5
$ws.Range("H136").Value = 5872.5
$ws.Range("I136").Value = 5872.5
$ws.Range("K136").Value = 17617.5
$ws.Range("M136").Value = -12517.5
$ws.Range("H137").Value = 5350
$ws.Range("I137").Value = 5937.5
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 17812.5
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -12712.5
$ws.Range("N137").Value = -19200
$ws.Range("H141").Value = 308598
$ws.Range("I141").Value = 753995
$ws.Range("J141").Value = 11666.667
$ws.Range("K141").Value = 2261985
$ws.Range("L141").Value = 35000.001
$ws.Range("M141").Value = -2256805
$ws.Range("N141").Value = -45360.001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 125002270
$ws.Range("I132").Value = 125002270
$ws.Range("K132").Value = 375006810
$ws.Range("M132").Value = -375004280

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 71666.336
$ws.Range("J36").Value = 71666.336
$ws.Range("L36").Value = 71666.336
$ws.Range("N36").Value = -72790.336
$ws.Range("H110").Value = 109817
$ws.Range("J110").Value = 109817
$ws.Range("L110").Value = 109817
$ws.Range("N110").Value = -117997
$ws.Range("H132").Value = 53912.6
$ws.Range("I132").Value = 5062.3335
$ws.Range("K132").Value = 15187.0005
$ws.Range("M132").Value = -12657.0005
$ws.Range("H136").Value = 41265.09
$ws.Range("I136").Value = 6890.5
$ws.Range("J136").Value = 233762.8
$ws.Range("K136").Value = 20671.5
$ws.Range("L136").Value = 701288.3999999999
$ws.Range("M136").Value = -18121.5
$ws.Range("N136").Value = -706388.3999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1322.7142
$ws.Range("I100").Value = 1393
$ws.Range("J100").Value = 1065
$ws.Range("K100").Value = 2786
$ws.Range("L100").Value = 2130
$ws.Range("M100").Value = -2245
$ws.Range("N100").Value = -3212
$ws.Range("H126").Value = 3756.919
$ws.Range("I126").Value = 3939.9678
$ws.Range("K126").Value = 11819.9034
$ws.Range("M126").Value = -9349.903399999999
